# Add 2060 perturbation columns (Year, IMAGE, MESSAGE, MiniCAM, MERGE, Policy)
# to both the CH4 table (rows 1-13) and the N2O table (rows 15-27), in the
# new column block AQ:AV, mirroring the existing 2010/2020/.../2300 blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535
$xlRight = -4152

# --- CH4 table (top block) ------------------------------------------------

# Row 2: block title
$ws.Cells.Item(2, 43).Value = "PAGE input: CH4 Shock, 2060"

# Row 3: sub-headers (Year, IMAGE, MESSAGE, MiniCAM, MERGE, Policy), right
# aligned on a yellow fill - same look as the other blocks' header row.
$ch4SubHeaders = @("Year", "IMAGE", "MESSAGE", "MiniCAM", "MERGE", "Policy")
for ($j = 0; $j -lt $ch4SubHeaders.Length; $j++) {
    $cell = $ws.Cells.Item(3, 43 + $j)
    $cell.Value = $ch4SubHeaders[$j]
    $cell.Interior.Color = $yellow
    $cell.HorizontalAlignment = $xlRight
}

# Rows 4-13: data (Year, then 5 model columns), each on a yellow fill.
$ch4Data = @(
    @(4, 2010, 0, 0, 0, 0, 0),
    @(5, 2020, 0, 0, 0, 0, 0),
    @(6, 2030, 0, 0, 0, 0, 0),
    @(7, 2040, 0, 0, 0, 0, 0),
    @(8, 2050, 0, 0, 0, 0, 0),
    @(9, 2060, 0.000071155643087639217, 0.000072944555823906132, 0.00007508513712286557, 0.000076260473212275143, 0.000088229397874300955),
    @(10, 2080, 0.000013418894915573354, 0.000014055531157619593, 0.000014544777158365107, 0.000014165668575349955, 0.000017364025824845309),
    @(11, 2100, 0.00000054923995536082519, 0.00000059687074613901105, 0.00000061714408356783323, 0.00000057665128834427647, 0.00000073930841433478277),
    @(12, 2200, 0.000000000090322458401459472, 0.000000000099343540060914399, 0.00000000010258888849534743, 0.000000000094655263627174685, 0.00000000012292778128752957),
    @(13, 2300, 0.00000000000012501111257279263, 0.00000000000013788969965844444, 0.00000000000014233059175694507, 0.00000000000013145040611561853, 0.00000000000017053025658242404)
)

foreach ($row in $ch4Data) {
    $r = $row[0]
    for ($j = 1; $j -lt $row.Length; $j++) {
        $cell = $ws.Cells.Item($r, 42 + $j)
        $cell.Value = $row[$j]
        $cell.Interior.Color = $yellow
    }
}

# --- N2O table (bottom block) ---------------------------------------------

# Row 16: block title
$ws.Cells.Item(16, 43).Value = "PAGE input: N2O Shock, 2060"

# Row 17: sub-headers
$n2oSubHeaders = @("Year", "IMAGE", "MESSAGE", "MiniCAM", "MERGE", "Policy")
for ($j = 0; $j -lt $n2oSubHeaders.Length; $j++) {
    $cell = $ws.Cells.Item(17, 43 + $j)
    $cell.Value = $n2oSubHeaders[$j]
    $cell.Interior.Color = $yellow
    $cell.HorizontalAlignment = $xlRight
}

# Rows 18-27: data
$n2oData = @(
    @(18, 2010, 0, 0, 0, 0, 0),
    @(19, 2020, 0, 0, 0, 0, 0),
    @(20, 2030, 0, 0, 0, 0, 0),
    @(21, 2040, 0, 0, 0, 0, 0),
    @(22, 2050, 0, 0, 0, 0, 0),
    @(23, 2060, 0.00034756994042368374, 0.00032965924978970161, 0.00032867020376360357, 0.00033633103433914903, 0.00034154397442547058),
    @(24, 2080, 0.00030442179878799878, 0.00028776869401197781, 0.00028413111429783646, 0.00029299055979585542, 0.00029949878678081908),
    @(25, 2100, 0.00018109292692546631, 0.00017303499667413347, 0.00016545806484028424, 0.00017257920543678807, 0.00017937031897590828),
    @(26, 2200, 0.000073388024830025605, 0.00007101522016244311, 0.000065696502927544692, 0.000069260370142663059, 0.000073215175532643782),
    @(27, 2300, 0.000045256296990203726, 0.000043949774506890638, 0.000040290512913099263, 0.000042597198992688767, 0.000045240913177369002)
)

foreach ($row in $n2oData) {
    $r = $row[0]
    for ($j = 1; $j -lt $row.Length; $j++) {
        $cell = $ws.Cells.Item($r, 42 + $j)
        $cell.Value = $row[$j]
        $cell.Interior.Color = $yellow
    }
}

# Keep the active selection where the new block starts (AR4), matching the
# saved workbook's view state.
$ws.Range("AR4").Select() | Out-Null

Write-Output "2060 CH4/N2O perturbation columns added"
